# closure_spreadsheet_main: add two new columns (Y: res_data_sent, Z: res_data_received)
# Y gets a header + per-row 0/1 flags for rows 2-43 (except row 22, left blank);
# Z only gets the header (no data rows yet) - matches the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Headers: copy the formatting of the existing header cell (X1) onto the
# two new header cells, then set their text. ---
$ws.Range("X1").Copy()
$ws.Range("Y1:Z1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("Y1").Value = "res_data_sent"
$ws.Range("Z1").Value = "res_data_received"

# --- Data column Y: 0/1 flag per university row (row 22 intentionally has no value) ---
$yValues = @{
    2=1; 3=1; 4=1; 5=0; 6=1; 7=1; 8=1; 9=0; 10=1; 11=1; 12=1; 13=1; 14=1; 15=1; 16=0;
    17=0; 18=1; 19=0; 20=1; 21=1; 23=0; 24=1; 25=1; 26=0; 27=1; 28=0; 29=1; 30=1; 31=1;
    32=1; 33=1; 34=1; 35=1; 36=1; 37=1; 38=1; 39=1; 40=1; 41=1; 42=1; 43=1
}

foreach ($row in $yValues.Keys) {
    $ws.Range("Y$row").Value = $yValues[$row]
}

# --- View state: keep column A frozen, and finish with the same selection as
# the authored workbook (Y38 in the scrolled/frozen right pane). ---
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Y38").Select()
